$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Liga de Fútbol de Primera División"
$ws.Range("B2").Value = "CS Herediano"
$ws.Range("D2").Value = "CS Cartagines"
$ws.Range("F2").Value = "43'"
$ws.Range("G2").Value = 43
$ws.Range("J2").Value = "23:46:01"

# Row 3
$ws.Range("A3").Value = "USL Championship"
$ws.Range("B3").Value = "Sacramento Republic FC"
$ws.Range("D3").Value = "Louisville City FC"
$ws.Range("F3").Value = "35'"
$ws.Range("G3").Value = 35
$ws.Range("J3").Value = "23:46:01"

# Row 4
$ws.Range("F4").Value = "80'"
$ws.Range("G4").Value = 80
$ws.Range("J4").Value = "23:46:02"
